# Auto-generated script applying the Famfrit_Profits data refresh diff.
# For each affected sheet, update the listed cells with their new values
# (market-board price/profit columns H, I, J, K, L, M, N refreshed by the scheduled runner).
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19: H19=4316.1665, J19=4579.4, L19=4579.4, N19=-4929.4
$ws.Range("H19").Value = 4316.1665
$ws.Range("J19").Value = 4579.4
$ws.Range("L19").Value = 4579.4
$ws.Range("N19").Value = -4929.4
# Row 40: H40=457214.1, I40=528722.5, J40=4327.3335, K40=528722.5, L40=4327.3335, M40=-528547.5, N40=-4677.3335
$ws.Range("H40").Value = 457214.1
$ws.Range("I40").Value = 528722.5
$ws.Range("J40").Value = 4327.3335
$ws.Range("K40").Value = 528722.5
$ws.Range("L40").Value = 4327.3335
$ws.Range("M40").Value = -528547.5
$ws.Range("N40").Value = -4677.3335
# Row 106: H106=2768.7144, I106=2676.2, K106=2676.2, M106=-2045.2
$ws.Range("H106").Value = 2768.7144
$ws.Range("I106").Value = 2676.2
$ws.Range("K106").Value = 2676.2
$ws.Range("M106").Value = -2045.2
# Row 132: H132=3549.4412, I132=3739.9333, K132=11219.7999, M132=-8689.7999
$ws.Range("H132").Value = 3549.4412
$ws.Range("I132").Value = 3739.9333
$ws.Range("K132").Value = 11219.7999
$ws.Range("M132").Value = -8689.7999
# Row 135: H135=916.5, I135=674.25, J135=1158.75, K135=6068.25, L135=10428.75, M135=-3533.25, N135=-15498.75
$ws.Range("H135").Value = 916.5
$ws.Range("I135").Value = 674.25
$ws.Range("J135").Value = 1158.75
$ws.Range("K135").Value = 6068.25
$ws.Range("L135").Value = 10428.75
$ws.Range("M135").Value = -3533.25
$ws.Range("N135").Value = -15498.75
# Row 138: H138=14927909, J138=25003606, L138=75010818, N138=-75021098
$ws.Range("H138").Value = 14927909
$ws.Range("J138").Value = 25003606
$ws.Range("L138").Value = 75010818
$ws.Range("N138").Value = -75021098

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 31: H31=63240.363, I31=12407.3, K31=12407.3, M31=-12113.3
$ws.Range("H31").Value = 63240.363
$ws.Range("I31").Value = 12407.3
$ws.Range("K31").Value = 12407.3
$ws.Range("M31").Value = -12113.3
# Row 32: H32=2555.952, I32=2496.8125, J32=4133, K32=2496.8125, L32=4133, M32=-2209.8125, N32=-4707
$ws.Range("H32").Value = 2555.952
$ws.Range("I32").Value = 2496.8125
$ws.Range("J32").Value = 4133
$ws.Range("K32").Value = 2496.8125
$ws.Range("L32").Value = 4133
$ws.Range("M32").Value = -2209.8125
$ws.Range("N32").Value = -4707
# Row 45: H45=4028.5, I45=3442, J45=5201.5, K45=3442, L45=5201.5, M45=-3065, N45=-5955.5
$ws.Range("H45").Value = 4028.5
$ws.Range("I45").Value = 3442
$ws.Range("J45").Value = 5201.5
$ws.Range("K45").Value = 3442
$ws.Range("L45").Value = 5201.5
$ws.Range("M45").Value = -3065
$ws.Range("N45").Value = -5955.5
# Row 61: H61=3787.25, I61=3567.25, K61=3567.25, M61=-3355.25
$ws.Range("H61").Value = 3787.25
$ws.Range("I61").Value = 3567.25
$ws.Range("K61").Value = 3567.25
$ws.Range("M61").Value = -3355.25
# Row 122: H122=2805.0454, I122=2458.5264, K122=7375.5792, M122=-4925.5792
$ws.Range("H122").Value = 2805.0454
$ws.Range("I122").Value = 2458.5264
$ws.Range("K122").Value = 7375.5792
$ws.Range("M122").Value = -4925.5792
# Row 136: H136=3787.25, I136=3567.25, K136=10701.75, M136=-8151.75
$ws.Range("H136").Value = 3787.25
$ws.Range("I136").Value = 3567.25
$ws.Range("K136").Value = 10701.75
$ws.Range("M136").Value = -8151.75

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 26: H26=13500.2, I26=13500.2, K26=13500.2, M26=-13208.2
$ws.Range("H26").Value = 13500.2
$ws.Range("I26").Value = 13500.2
$ws.Range("K26").Value = 13500.2
$ws.Range("M26").Value = -13208.2
# Row 94: H94=2431, I94=2431, K94=2431, M94=-1980
$ws.Range("H94").Value = 2431
$ws.Range("I94").Value = 2431
$ws.Range("K94").Value = 2431
$ws.Range("M94").Value = -1980
# Row 99: H99=5463.857, I99=3976.6667, J99=6579.25, K99=3976.6667, L99=6579.25, M99=-2478.6667, N99=-9575.25
$ws.Range("H99").Value = 5463.857
$ws.Range("I99").Value = 3976.6667
$ws.Range("J99").Value = 6579.25
$ws.Range("K99").Value = 3976.6667
$ws.Range("L99").Value = 6579.25
$ws.Range("M99").Value = -2478.6667
$ws.Range("N99").Value = -9575.25
# Row 105: H105=10033.814, J105=7389.706, L105=7389.706, N105=-10883.706
$ws.Range("H105").Value = 10033.814
$ws.Range("J105").Value = 7389.706
$ws.Range("L105").Value = 7389.706
$ws.Range("N105").Value = -10883.706
# Row 107: H107=3178.8572, I107=2910.5454, J107=4162.6665, K107=2910.5454, L107=4162.6665, M107=-990.5454, N107=-8002.6665
$ws.Range("H107").Value = 3178.8572
$ws.Range("I107").Value = 2910.5454
$ws.Range("J107").Value = 4162.6665
$ws.Range("K107").Value = 2910.5454
$ws.Range("L107").Value = 4162.6665
$ws.Range("M107").Value = -990.5454
$ws.Range("N107").Value = -8002.6665

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 6: H6=1001, I6=1001, K6=1001, M6=-888
$ws.Range("H6").Value = 1001
$ws.Range("I6").Value = 1001
$ws.Range("K6").Value = 1001
$ws.Range("M6").Value = -888
# Row 16: H16=1218.2, I16=1271.125, K16=1271.125, M16=-984.125
$ws.Range("H16").Value = 1218.2
$ws.Range("I16").Value = 1271.125
$ws.Range("K16").Value = 1271.125
$ws.Range("M16").Value = -984.125
# Row 32: H32=1450, J32=1450, L32=1450, N32=-2082
$ws.Range("H32").Value = 1450
$ws.Range("J32").Value = 1450
$ws.Range("L32").Value = 1450
$ws.Range("N32").Value = -2082
# Row 113: H113=1218.2, I113=1271.125, K113=1271.125, M113=898.875
$ws.Range("H113").Value = 1218.2
$ws.Range("I113").Value = 1271.125
$ws.Range("K113").Value = 1271.125
$ws.Range("M113").Value = 898.875
# Row 122: H122=1699.0714, I122=1565.8334, K122=4697.5002, M122=-2247.5002
$ws.Range("H122").Value = 1699.0714
$ws.Range("I122").Value = 1565.8334
$ws.Range("K122").Value = 4697.5002
$ws.Range("M122").Value = -2247.5002

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 70: H70=500, I70=500, K70=1500, M70=-1185
$ws.Range("H70").Value = 500
$ws.Range("I70").Value = 500
$ws.Range("K70").Value = 1500
$ws.Range("M70").Value = -1185
# Row 73: H73=500, I73=500, K73=1500, M73=-408
$ws.Range("H73").Value = 500
$ws.Range("I73").Value = 500
$ws.Range("K73").Value = 1500
$ws.Range("M73").Value = -408
# Row 75: H75=1086.6, I75=996, J75=1222.5, K75=2988, L75=3667.5, M75=-1990, N75=-5663.5
$ws.Range("H75").Value = 1086.6
$ws.Range("I75").Value = 996
$ws.Range("J75").Value = 1222.5
$ws.Range("K75").Value = 2988
$ws.Range("L75").Value = 3667.5
$ws.Range("M75").Value = -1990
$ws.Range("N75").Value = -5663.5
# Row 78: H78=1086.6, I78=996, J78=1222.5, K78=8964, L78=11002.5, M78=-3972, N78=-20986.5
$ws.Range("H78").Value = 1086.6
$ws.Range("I78").Value = 996
$ws.Range("J78").Value = 1222.5
$ws.Range("K78").Value = 8964
$ws.Range("L78").Value = 11002.5
$ws.Range("M78").Value = -3972
$ws.Range("N78").Value = -20986.5
# Row 132: H132=3906.1333, I132=1371.5, J132=5595.8887, K132=12343.5, L132=50362.99830000001, M132=-9813.5, N132=-55422.99830000001
$ws.Range("H132").Value = 3906.1333
$ws.Range("I132").Value = 1371.5
$ws.Range("J132").Value = 5595.8887
$ws.Range("K132").Value = 12343.5
$ws.Range("L132").Value = 50362.99830000001
$ws.Range("M132").Value = -9813.5
$ws.Range("N132").Value = -55422.99830000001

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80: H80=1454, J80=1604.8, L80=1604.8, N80=-3600.8
$ws.Range("H80").Value = 1454
$ws.Range("J80").Value = 1604.8
$ws.Range("L80").Value = 1604.8
$ws.Range("N80").Value = -3600.8
# Row 83: H83=1454, J83=1604.8, L83=8024, N83=-18008
$ws.Range("H83").Value = 1454
$ws.Range("J83").Value = 1604.8
$ws.Range("L83").Value = 8024
$ws.Range("N83").Value = -18008
# Row 97: H97=2502, I97=1999.5, J97=2837, K97=1999.5, L97=2837, M97=-1503.5, N97=-3829
$ws.Range("H97").Value = 2502
$ws.Range("I97").Value = 1999.5
$ws.Range("J97").Value = 2837
$ws.Range("K97").Value = 1999.5
$ws.Range("L97").Value = 2837
$ws.Range("M97").Value = -1503.5
$ws.Range("N97").Value = -3829
# Row 102: H102=5332.737, I102=2439.5454, K102=2439.5454, M102=-817.5454
$ws.Range("H102").Value = 5332.737
$ws.Range("I102").Value = 2439.5454
$ws.Range("K102").Value = 2439.5454
$ws.Range("M102").Value = -817.5454
# Row 113: H113=3008.3333, J113=4602.75, L113=4602.75, N113=-8942.75
$ws.Range("H113").Value = 3008.3333
$ws.Range("J113").Value = 4602.75
$ws.Range("L113").Value = 4602.75
$ws.Range("N113").Value = -8942.75
# Row 122: H122=1399.9642, I122=1029.5652, K122=3088.6956, M122=-638.6956
$ws.Range("H122").Value = 1399.9642
$ws.Range("I122").Value = 1029.5652
$ws.Range("K122").Value = 3088.6956
$ws.Range("M122").Value = -638.6956
# Row 123: H123=59498.5, J123=64997.332, L123=64997.332, N123=-69897.33199999999
$ws.Range("H123").Value = 59498.5
$ws.Range("J123").Value = 64997.332
$ws.Range("L123").Value = 64997.332
$ws.Range("N123").Value = -69897.33199999999

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 61: H61=21741822, I61=29414202, J61=3407.3333, K61=29414202, L61=3407.3333, M61=-29414000, N61=-3811.3333
$ws.Range("H61").Value = 21741822
$ws.Range("I61").Value = 29414202
$ws.Range("J61").Value = 3407.3333
$ws.Range("K61").Value = 29414202
$ws.Range("L61").Value = 3407.3333
$ws.Range("M61").Value = -29414000
$ws.Range("N61").Value = -3811.3333
# Row 68: H68=2750, I68=1250, J68=4250, K68=1250, L68=4250, M68=-501, N68=-5748
$ws.Range("H68").Value = 2750
$ws.Range("I68").Value = 1250
$ws.Range("J68").Value = 4250
$ws.Range("K68").Value = 1250
$ws.Range("L68").Value = 4250
$ws.Range("M68").Value = -501
$ws.Range("N68").Value = -5748
# Row 71: H71=2750, I71=1250, J71=4250, K71=6250, L71=21250, M71=-2506, N71=-28738
$ws.Range("H71").Value = 2750
$ws.Range("I71").Value = 1250
$ws.Range("J71").Value = 4250
$ws.Range("K71").Value = 6250
$ws.Range("L71").Value = 21250
$ws.Range("M71").Value = -2506
$ws.Range("N71").Value = -28738
# Row 82: H82=2723.6365, I82=2741.3076, J82=2698.111, K82=2741.3076, L82=2698.111, M82=-2380.3076, N82=-3420.111
$ws.Range("H82").Value = 2723.6365
$ws.Range("I82").Value = 2741.3076
$ws.Range("J82").Value = 2698.111
$ws.Range("K82").Value = 2741.3076
$ws.Range("L82").Value = 2698.111
$ws.Range("M82").Value = -2380.3076
$ws.Range("N82").Value = -3420.111
# Row 85: H85=2723.6365, I85=2741.3076, J85=2698.111, K85=2741.3076, L85=2698.111, M85=-1493.3076, N85=-5194.111
$ws.Range("H85").Value = 2723.6365
$ws.Range("I85").Value = 2741.3076
$ws.Range("J85").Value = 2698.111
$ws.Range("K85").Value = 2741.3076
$ws.Range("L85").Value = 2698.111
$ws.Range("M85").Value = -1493.3076
$ws.Range("N85").Value = -5194.111
# Row 93: H93=4017.5, I93=3922, K93=3922, M93=-2674
$ws.Range("H93").Value = 4017.5
$ws.Range("I93").Value = 3922
$ws.Range("K93").Value = 3922
$ws.Range("M93").Value = -2674
# Row 100: H100=3704.375, I100=3380.762, K100=3380.762, M100=-2839.762
$ws.Range("H100").Value = 3704.375
$ws.Range("I100").Value = 3380.762
$ws.Range("K100").Value = 3380.762
$ws.Range("M100").Value = -2839.762
# Row 113: H113=21741822, I113=29414202, J113=3407.3333, K113=29414202, L113=3407.3333, M113=-29412032, N113=-7747.3333
$ws.Range("H113").Value = 21741822
$ws.Range("I113").Value = 29414202
$ws.Range("J113").Value = 3407.3333
$ws.Range("K113").Value = 29414202
$ws.Range("L113").Value = 3407.3333
$ws.Range("M113").Value = -29412032
$ws.Range("N113").Value = -7747.3333
# Row 122: H122=4878.8667, I122=4312.091, K122=12936.273, M122=-10486.273
$ws.Range("H122").Value = 4878.8667
$ws.Range("I122").Value = 4312.091
$ws.Range("K122").Value = 12936.273
$ws.Range("M122").Value = -10486.273
# Row 133: H133=45705, J133=45705, L133=45705, N133=-50765
$ws.Range("H133").Value = 45705
$ws.Range("J133").Value = 45705
$ws.Range("L133").Value = 45705
$ws.Range("N133").Value = -50765

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 40: H40=10539.857, I40=7755.8, J40=17500, K40=7755.8, L40=17500, M40=-7606.8, N40=-17798
$ws.Range("H40").Value = 10539.857
$ws.Range("I40").Value = 7755.8
$ws.Range("J40").Value = 17500
$ws.Range("K40").Value = 7755.8
$ws.Range("L40").Value = 17500
$ws.Range("M40").Value = -7606.8
$ws.Range("N40").Value = -17798
# Row 81: H81=1566.9286, J81=2087.2, L81=4174.4, N81=-6296.4
$ws.Range("H81").Value = 1566.9286
$ws.Range("J81").Value = 2087.2
$ws.Range("L81").Value = 4174.4
$ws.Range("N81").Value = -6296.4
# Row 84: H84=1566.9286, J84=2087.2, L84=20872, N84=-31480
$ws.Range("H84").Value = 1566.9286
$ws.Range("J84").Value = 2087.2
$ws.Range("L84").Value = 20872
$ws.Range("N84").Value = -31480
# Row 96: H96=3433.375, I96=2979.1538, K96=2979.1538, M96=-1606.1538
$ws.Range("H96").Value = 3433.375
$ws.Range("I96").Value = 2979.1538
$ws.Range("K96").Value = 2979.1538
$ws.Range("M96").Value = -1606.1538
# Row 113: H113=1178.238, I113=1125, K113=3375, M113=-1205
$ws.Range("H113").Value = 1178.238
$ws.Range("I113").Value = 1125
$ws.Range("K113").Value = 3375
$ws.Range("M113").Value = -1205
# Row 133: H133=84320.60000000001, J133=84320.60000000001, L133=84320.60000000001, N133=-94440.60000000001
$ws.Range("H133").Value = 84320.60000000001
$ws.Range("J133").Value = 84320.60000000001
$ws.Range("L133").Value = 84320.60000000001
$ws.Range("N133").Value = -94440.60000000001

